# Adds two new columns "I0" (I) and "IF" (J) to the sheet, matching the
# existing header style used by the other header cells (copy H1's
# formatting so the same style index is reused), and fills in the
# numeric data for rows 2..52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---
# Copy the formatting of an existing header cell (H1) onto the new
# header cells so they match (bold, centered, bordered) and then set
# their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data values (rows 2..52) ---
$iVals = @(9,4,9,9,7,7,11,8,7,7,6,2,6,8,9,6,7,5,10,6,10,7,8,5,6,5,9,4,9,7,6,6,6,8,7,6,6,6,7,5,5,7,6,6,8,7,6,8,4,7,5)
$jVals = @(9,6,9,9,7,8,11,8,7,7,6,3,6,8,9,6,7,6,10,8,10,7,8,6,7,6,9,5,9,7,6,7,7,9,7,6,6,7,7,6,5,7,7,7,9,7,7,8,5,7,5)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
